$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '41.752.12'
$ws.Range("E2").Value = '  +5.81%  '

$ws.Range("D3").Value = '2.263.53'
$ws.Range("E3").Value = '  +4.63%  '

$ws.Range("E4").Value = '  +0.03%  '

Set-TextValue $ws.Range("D5") '234.20'
$ws.Range("E5").Value = '  +2.44%  '

$ws.Range("E6").Value = '  +3.55%  '

Set-TextValue $ws.Range("D7") '64.44'
$ws.Range("E7").Value = '  +1.23%  '

$ws.Range("E9").Value = '  +4.31%  '

Set-TextValue $ws.Range("D10") '60.11'
$ws.Range("E10").Value = '  +3.47%  '

Set-TextValue $ws.Range("D11") '0.0898'
$ws.Range("E11").Value = '  +5.14%  '

$ws.Range("E12").Value = '  +2.40%  '

$ws.Range("D13").Value = '2.601.60'
$ws.Range("E13").Value = '  +4.67%  '

Set-TextValue $ws.Range("D14") '16.16'
$ws.Range("E14").Value = '  +0.53%  '

Set-TextValue $ws.Range("D15") '22.80'
$ws.Range("E15").Value = '  +3.24%  '

Set-TextValue $ws.Range("D16") '0.829'
$ws.Range("E16").Value = '  +1.93%  '

Set-TextValue $ws.Range("D17") '5.70'
$ws.Range("E17").Value = '  +3.21%  '

$ws.Range("D18").Value = '2.264.89'
$ws.Range("E18").Value = '  +4.23%  '

$ws.Range("D19").Value = '41.626.63'
$ws.Range("E19").Value = '  +5.50%  '

$ws.Range("D20").Value = '0.0₃0941'
$ws.Range("E20").Value = '  +10.85%  '

Set-TextValue $ws.Range("D21") '74.95'
$ws.Range("E21").Value = '  +4.29%  '

Set-TextValue $ws.Range("D22") '6.20'
$ws.Range("E22").Value = '  -0.28%  '

Set-TextValue $ws.Range("D23") '252.25'
$ws.Range("E23").Value = '  +9.94%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E25").Value = '  +3.14%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D26") '2.35'
$ws.Range("E26").Value = '  +0.42%  '

Set-TextValue $ws.Range("D27") '9.88'
$ws.Range("E27").Value = '  +3.88%  '

$ws.Range("E28").Value = '  +5.43%  '

Set-TextValue $ws.Range("D29") '171.30'
$ws.Range("E29").Value = '  -0.50%  '

Set-TextValue $ws.Range("D30") '20.58'
$ws.Range("E30").Value = '  +3.56%  '

$ws.Range("E31").Value = '  +1.80%  '

Set-TextValue $ws.Range("D32") '2.84'
$ws.Range("E32").Value = '  +5.98%  '

$ws.Range("E33").Value = '  +2.98%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D34") '5.11'
$ws.Range("E34").Value = '  +8.37%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D35") '4.79'
$ws.Range("E35").Value = '  +3.62%  '

$ws.Range("E36").Value = '  +3.20%  '

Set-TextValue $ws.Range("D37") '6.93'
$ws.Range("E37").Value = '  -2.22%  '

Set-TextValue $ws.Range("D38") '3.88'
$ws.Range("E38").Value = '  +8.20%  '

Set-TextValue $ws.Range("D39") '2.48'
$ws.Range("E39").Value = '  +1.54%  '

$ws.Range("E40").Value = '  +60.13%  '

Set-TextValue $ws.Range("D41") '5.16'
$ws.Range("E41").Value = '  +20.58%  '

$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("E43").Value = '  +5.44%  '

$ws.Range("E44").Value = '  +13.21%  '

Set-TextValue $ws.Range("D45") '102.87'
$ws.Range("E45").Value = '  -0.29%  '

Set-TextValue $ws.Range("D46") '0.0991'
$ws.Range("E46").Value = '  +6.99%  '

Set-TextValue $ws.Range("D47") '17.71'
$ws.Range("E47").Value = '  -0.64%  '

$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("D49").Value = '1.510.07'
$ws.Range("E49").Value = '  -0.89%  '

$ws.Range("E50").Value = '  +2.36%  '

Set-TextValue $ws.Range("D51") '2.79'
$ws.Range("E51").Value = '  -1.09%  '
